$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main_page")

$ws.Range("D4").Value = 'text = "Active in last 24 hours", colour = "purple", icon = "clock"'
$ws.Range("D5").Value = 'text = "Active in last 7 days", colour = "green", icon = "calendar"'
$ws.Range("D2").Value = 'text = "Total", colour = "yellow", icon = "user"'
$ws.Range("D3").Value = 'text = "Consented", colour = "aqua", icon = "clipboard"'

$ws.Range("D4").Select()
